# Restore revision change: Rules sheet, row "R40" (row 10), the "From" value
# (cell C10) reverts from 18 back to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
